# Update countries & provincias Spain
# - Re-rank a few countries (Panama above Republica Dominicana,
#   Nueva Zelanda above Kazajistan, Niger above Republica de Yibuti)
#   by swapping their row contents, and refresh case totals for the
#   affected rows plus a couple of standalone updates (Mexico,
#   Islas Caimanes, Zimbabue).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($row, $pais, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($row, 1).Value = $pais
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Mexico: refreshed totals, same rank
Set-Row 38 "Mexico" 6297 450 2125 3686 207 37 486

# Panama overtakes Republica Dominicana
Set-Row 48 "Panama" 4016 265 98 3809 99 6 109
Set-Row 49 "Republica Dominicana" 3755 0 215 3344 121 0 196

# Nueva Zelanda overtakes Kazajistan
Set-Row 69 "Nueva Zelanda" 1409 8 816 582 2 2 11
Set-Row 70 "Kazajistan" 1402 0 277 1108 22 0 17

# Niger overtakes Republica de Yibuti
Set-Row 94 "Niger" 609 25 105 489 0 1 15
Set-Row 95 "Republica de Yibuti" 591 0 73 516 0 0 2

# Islas Caimanes: refreshed totals, same rank
Set-Row 148 "Islas Caimanes" 61 1 7 53 3 0 1

# Zimbabue: refreshed totals, same rank
Set-Row 173 "Zimbabue" 24 1 2 19 0 0 3
